$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, copying the header style (s="1") from E1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# time_taken values for rows 2-82 (one per gene/panel row)
$timeTaken = @(
    "2021-10-05 10:50:49.330704",
    "2021-10-05 10:50:49.330716",
    "2021-10-05 10:50:49.330719",
    "2021-10-05 10:50:49.330722",
    "2021-10-05 10:50:49.330725",
    "2021-10-05 10:50:49.330728",
    "2021-10-05 10:50:49.330730",
    "2021-10-05 10:50:49.330733",
    "2021-10-05 10:50:49.330735",
    "2021-10-05 10:50:49.330738",
    "2021-10-05 10:50:49.330740",
    "2021-10-05 10:50:49.330743",
    "2021-10-05 10:50:49.330745",
    "2021-10-05 10:50:49.330748",
    "2021-10-05 10:50:49.330750",
    "2021-10-05 10:50:49.330753",
    "2021-10-05 10:50:49.330755",
    "2021-10-05 10:50:49.330758",
    "2021-10-05 10:50:49.330761",
    "2021-10-05 10:50:49.330763",
    "2021-10-05 10:50:49.330765",
    "2021-10-05 10:50:49.330768",
    "2021-10-05 10:50:49.330770",
    "2021-10-05 10:50:49.330773",
    "2021-10-05 10:50:49.330775",
    "2021-10-05 10:50:49.330778",
    "2021-10-05 10:50:49.330780",
    "2021-10-05 10:50:49.330783",
    "2021-10-05 10:50:49.330785",
    "2021-10-05 10:50:49.330788",
    "2021-10-05 10:50:49.330790",
    "2021-10-05 10:50:49.330793",
    "2021-10-05 10:50:49.330796",
    "2021-10-05 10:50:49.330798",
    "2021-10-05 10:50:49.330801",
    "2021-10-05 10:50:49.330803",
    "2021-10-05 10:50:49.330806",
    "2021-10-05 10:50:49.330808",
    "2021-10-05 10:50:49.330811",
    "2021-10-05 10:50:49.330813",
    "2021-10-05 10:50:49.330816",
    "2021-10-05 10:50:49.330819",
    "2021-10-05 10:50:49.330821",
    "2021-10-05 10:50:49.330824",
    "2021-10-05 10:50:49.330826",
    "2021-10-05 10:50:49.330829",
    "2021-10-05 10:50:49.330831",
    "2021-10-05 10:50:49.330834",
    "2021-10-05 10:50:49.330836",
    "2021-10-05 10:50:49.330839",
    "2021-10-05 10:50:49.330841",
    "2021-10-05 10:50:49.330844",
    "2021-10-05 10:50:49.330846",
    "2021-10-05 10:50:49.330849",
    "2021-10-05 10:50:49.330851",
    "2021-10-05 10:50:49.330854",
    "2021-10-05 10:50:49.330856",
    "2021-10-05 10:50:49.330859",
    "2021-10-05 10:50:49.330861",
    "2021-10-05 10:50:49.330864",
    "2021-10-05 10:50:49.330866",
    "2021-10-05 10:50:49.330869",
    "2021-10-05 10:50:49.330871",
    "2021-10-05 10:50:49.330874",
    "2021-10-05 10:50:49.330877",
    "2021-10-05 10:50:49.330880",
    "2021-10-05 10:50:49.330882",
    "2021-10-05 10:50:49.330885",
    "2021-10-05 10:50:49.330887",
    "2021-10-05 10:50:49.330890",
    "2021-10-05 10:50:49.330892",
    "2021-10-05 10:50:49.330895",
    "2021-10-05 10:50:49.330897",
    "2021-10-05 10:50:49.330900",
    "2021-10-05 10:50:49.330902",
    "2021-10-05 10:50:49.330905",
    "2021-10-05 10:50:49.330909",
    "2021-10-05 10:50:49.330912",
    "2021-10-05 10:50:49.330915",
    "2021-10-05 10:50:49.330917",
    "2021-10-05 10:50:49.330920"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
